# Weekly driver report update for 2025-05-05
# Updates the "Good Drivers (Roaming > 99.8%)" table (rows 13-24) on the
# "Driver Summary" sheet with the latest weekly roaming-impact numbers.
# Some adapter/driver rows were re-ordered by the upstream report as part
# of this refresh, so both the text/label columns and the numeric columns
# are rewritten per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.0.4
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.0.4"
$ws.Range("B13").Value = 1293197
$ws.Range("C13").Value = 4322
$ws.Range("D13").Value = 1009
$ws.Range("E13").Value = 1990
$ws.Range("F13").Value = 1298528
$ws.Range("H13").Value = "22.250.0.4"
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = "2023-07-25"

# Row 14 - Intel(R) Wi-Fi 6E AX211 160MHz - 22.220.0.4
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.220.0.4"
$ws.Range("B14").Value = 31517
$ws.Range("C14").Value = 112
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 66
$ws.Range("F14").Value = 31629
$ws.Range("H14").Value = "22.220.0.4"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "2023-03-28"

# Row 15 - Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8"
$ws.Range("B15").Value = 467311
$ws.Range("C15").Value = 772
$ws.Range("D15").Value = 931
$ws.Range("E15").Value = 706
$ws.Range("F15").Value = 469014
$ws.Range("H15").Value = "23.10.0.8"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "2023-10-30"

# Row 16 - Intel(R) Wi-Fi 6E AX211 160MHz - 23.120.0.3
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.120.0.3"
$ws.Range("B16").Value = 455081
$ws.Range("C16").Value = 1861
$ws.Range("D16").Value = 52
$ws.Range("E16").Value = 639
$ws.Range("F16").Value = 456994
$ws.Range("H16").Value = "23.120.0.3"
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "2025-02-05"

# Row 17 - Intel(R) Wi-Fi 6E AX211 160MHz - 22.230.0.8
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.230.0.8"
$ws.Range("B17").Value = 1787924
$ws.Range("C17").Value = 3326
$ws.Range("D17").Value = 2614
$ws.Range("E17").Value = 3038
$ws.Range("F17").Value = 1793864
$ws.Range("H17").Value = "22.230.0.8"
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "2023-05-08"

# Row 18 - Intel(R) Wi-Fi 6E AX211 160MHz - 23.70.2.3
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.70.2.3"
$ws.Range("B18").Value = 218767
$ws.Range("C18").Value = 334
$ws.Range("D18").Value = 313
$ws.Range("E18").Value = 573
$ws.Range("F18").Value = 219414
$ws.Range("H18").Value = "23.70.2.3"
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "2024-07-23"

# Row 19 - Intel(R) Wi-Fi 6E AX211 160MHz - 22.110.1.1
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.110.1.1"
$ws.Range("B19").Value = 135467
$ws.Range("C19").Value = 189
$ws.Range("D19").Value = 263
$ws.Range("E19").Value = 196
$ws.Range("F19").Value = 135919
$ws.Range("H19").Value = "22.110.1.1"
$ws.Range("J19").NumberFormat = "@"
$ws.Range("J19").Value = "2022-01-01"

# Row 20 - Intel(R) Wi-Fi 6E AX211 160MHz - 23.100.0.4
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.100.0.4"
$ws.Range("B20").Value = 240434
$ws.Range("C20").Value = 421
$ws.Range("D20").Value = 37
$ws.Range("E20").Value = 409
$ws.Range("F20").Value = 240892
$ws.Range("H20").Value = "23.100.0.4"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "2024-11-10"

# Row 21 - Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.1.3
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.1.3"
$ws.Range("B21").Value = 151287
$ws.Range("C21").Value = 285
$ws.Range("D21").Value = 75
$ws.Range("E21").Value = 332
$ws.Range("F21").Value = 151647
$ws.Range("H21").Value = "23.80.1.3"
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = "2024-09-03"

# Row 22 - Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1
$ws.Range("A22").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B22").Value = 272039
$ws.Range("C22").Value = 213
$ws.Range("D22").Value = 131
$ws.Range("E22").Value = 316
$ws.Range("F22").Value = 272383
$ws.Range("H22").Value = "22.100.1.1"
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value = "2022-05-01"

# Row 23 - Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3
$ws.Range("A23").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B23").Value = 14561
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 59
$ws.Range("F23").Value = 14561
$ws.Range("H23").Value = "22.150.0.3"
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = "2022-05-23"

# Row 24 - Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1
$ws.Range("A24").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B24").Value = 12018
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 61
$ws.Range("F24").Value = 12018
$ws.Range("H24").Value = "22.150.3.1"
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value = "2022-08-29"
